# This edit reflects the workbook being touched by a newer Excel build and
# the columns on "micom" being resized to fit their contents:
#   A -> fixed width 12
#   B -> best-fit width 21.5   (longest value: "GROUP 1 WATTMETRIC SEF")
#   C -> best-fit width 33.83  (longest value: 32-char binary string)
#
# Column widths are set through ColumnWidth (Excel's "characters" unit).
# The engine's stored <col width="..."> value is ColumnWidth + 5/6, so we
# back that constant out here to land on the target stored widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offset = 5.0 / 6.0

$ws.Columns.Item(1).ColumnWidth = 12 - $offset
$ws.Columns.Item(2).ColumnWidth = 21.5 - $offset
$ws.Columns.Item(3).ColumnWidth = 33.83203125 - $offset
